$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 42/43 also swap positions (Monero <-> EnergySwap) in addition to value updates.
# D-column price values are prefixed with a leading apostrophe so Excel keeps them
# as literal text (matching the source data) instead of re-parsing numeric-looking
# strings like "319.73" into a Number cell.

$ws.Range("D2").Value = "'47.108.83"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "'2.478.82"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'319.73"
$ws.Range("E5").Value = "  -1.42%  "
$ws.Range("D6").Value = "'107.91"
$ws.Range("E6").Value = "  +2.49%  "
$ws.Range("D7").Value = "'0.521"
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.532"
$ws.Range("E9").Value = "  -1.98%  "
$ws.Range("D10").Value = "'38.72"
$ws.Range("E10").Value = "  +6.75%  "
$ws.Range("E11").Value = "  -1.34%  "
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").Value = "'18.12"
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("E14").Value = "  -0.63%  "
$ws.Range("D15").Value = "'2.845.31"
$ws.Range("E15").Value = "  -0.76%  "
$ws.Range("D16").Value = "'2.474.78"
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("E17").Value = "  -0.31%  "
$ws.Range("D18").Value = "'47.038.44"
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("E20").Value = "  +1.69%  "
$ws.Range("D21").Value = "'0.0₃0929"
$ws.Range("E21").Value = "  -0.92%  "
$ws.Range("E22").Value = "  +13.49%  "
$ws.Range("D23").Value = "'70.16"
$ws.Range("E23").Value = "  -0.66%  "
$ws.Range("D24").Value = "'244.32"
$ws.Range("E24").Value = "  -2.71%  "
$ws.Range("E25").Value = "  -0.38%  "
$ws.Range("D27").Value = "'25.59"
$ws.Range("E27").Value = "  -2.63%  "
$ws.Range("E28").Value = "  +3.25%  "
$ws.Range("D29").Value = "'9.99"
$ws.Range("E29").Value = "  +1.36%  "
$ws.Range("D30").Value = "'34.81"
$ws.Range("E30").Value = "  -0.72%  "
$ws.Range("D31").Value = "'0.133"
$ws.Range("E31").Value = "  -1.96%  "
$ws.Range("D32").Value = "'49.41"
$ws.Range("E32").Value = "  -0.40%  "
$ws.Range("D33").Value = "'19.89"
$ws.Range("E33").Value = "  +0.85%  "
$ws.Range("E34").Value = "  +0.32%  "
$ws.Range("D35").Value = "'0.0780"
$ws.Range("E35").Value = "  +1.19%  "
$ws.Range("E36").Value = "  +0.25%  "
$ws.Range("E37").Value = "  +1.86%  "
$ws.Range("D38").Value = "'4.62"
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("E39").Value = "  -0.98%  "
$ws.Range("E40").Value = "  -0.35%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "'119.35"
$ws.Range("E42").Value = "  -2.87%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'21.79"
$ws.Range("E43").Value = "  +3.68%  "
$ws.Range("D44").Value = "'0.0293"
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("D45").Value = "'1.976.36"
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("E46").Value = "  +0.63%  "
$ws.Range("E47").Value = "  -4.58%  "
$ws.Range("D48").Value = "'9.06"
$ws.Range("E48").Value = "  +0.45%  "
$ws.Range("E49").Value = "  -2.42%  "
$ws.Range("D50").Value = "'5.12"
$ws.Range("E50").Value = "  -4.82%  "
$ws.Range("D51").Value = "'57.15"
$ws.Range("E51").Value = "  +4.77%  "
